$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A252").Value = "2023-12-12 19:13:13"
$ws.Range("B252").Value = 0.002

$ws.Range("A253").Value = "2023-12-12 19:13:53"
$ws.Range("B253").Value = 0.0024

$ws.Range("A254").Value = "2023-12-12 19:14:05"
$ws.Range("B254").Value = 0.0006000000000000001
